# Add "NA" values in the duplicate_image_filename column (E) for all
# data rows (rows 2-21) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E21").Value = "NA"
